$d = $word.ActiveDocument

function Assert-True($cond, $msg) {
    if (-not $cond) { throw $msg }
}

# --- Paragraph 1: header/title line, two runs split by <w:br/> ---
# These runs contain emoji outside the BMP. Range.Start/End use UTF-16 code-unit
# offsets (standard Word COM convention and consistent with .NET string .Length),
# so build the two sub-ranges from precomputed (.NET) string lengths and assign
# .Text directly -- this also sidesteps the Find-replace auto "smart quote" pass.
$old1 = '🚀המאמר היומי של מייק 17.09.24: ⚡️🚀'
$new1 = '⚡️🚀המאמר היומי של מייק 16.09.24: ⚡️🚀'
$old2 = 'STaR: Self-Taught Reasoner Bootstrapping Reasoning With Reasoning'
$new2 = 'Rethinking Benchmark and Contamination for Language Models with Rephrased Samples'

$r1 = $d.Range(0, $old1.Length)
Assert-True ($r1.Text -eq $old1) "paragraph 1 run 1 text mismatch"
$r1.Text = $new1

$r2start = $new1.Length + 1   # +1 skips the <w:br/> line break
$r2end = $r2start + $old2.Length
$r2 = $d.Range($r2start, $r2end)
Assert-True ($r2.Text -eq $old2) "paragraph 1 run 2 text mismatch"
$r2.Text = $new2

# --- Paragraph 2 ---
$p2 = $d.Paragraphs.Item(2).Range
$found2 = $p2.Find.Execute('אני ממשיך לחפור במאמרי שאולי עיצבו את הנתיב הובילו ל-o1 של openai. הפעם נברתי כה עמוק שהגעתי למאמר שיצא לפני שנתיים וחצי (בדיפ היום זה כמו 100 שנה במתמטיקה). שימו לב שהמאמר יצא עוד לפני chatgpt. המאמר הזה מציע שיטה לשיפור יכולת reasoning של מודל שפה כאשר בידנו יש דאטהסט גדול של שאלות ותשובות D ודאטהסט קטן D_R הרבה יותר (המאמר מדבר על 10 דוגמאות בלבד) המכיל בנוסף גם את שרשרת ה-reasoning.', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Assert-True $found2 "paragraph 2 old text not found"
$p2.Text = 'חתיכת נושא זה. לאחרונה אני ניהלתי מספר שיחות עם אנשי NLP לא מעטים על הנושא הזה. מי שעוקב אחריי ברשתות החברתיות אולי שם לב כי אני בד״כ לא מתלהב ממודל שפה שניצח את כל המודלים הקיימים בכל הבנצ''מרקים. הסיבה לכך היא די טבעית ונובעת מכך שבלא מעט מקרים לא מפרסמים באופן גלוי את כל הדאטה שעליה המודל אומן. '

# --- Paragraph 3 ---
$p3 = $d.Paragraphs.Item(3).Range
$found3 = $p3.Find.Execute('כאשר אני מדבר על שיפור איכות ה-reasoning אני בעצם מתכוון לפיינטיון של המודל במטרה לקבל מודל חזק יותר ב-reasoning. המחברים מציעים אלגוריתם המורכב משני שלבים עיקריים. בשלב הראשון מזינים את הבאץ'' של שאלות למודל שפה כאשר בנוסף לשאלות הפרומפט מכיל את דוגמאות לשרשראות ה-reasoning m מ- D_R. המודל מתבקש לבנות שרשרת reasoning לכל השאלות מבאץ'' (לא מ-D_R) ולהגיע לתשובה הסופית.', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Assert-True $found3 "paragraph 3 old text not found"
$p3.Text = 'כמובן שהחשד שלי הוא הדאטה(משימות) האימון יהיו דומות מדי לאלו שמופיעות בבנצ''מרקים האלו. כמובן אני לא בא להאשים אנשים על כך שהם מרמים בכוונה (למרות שבטח יש מקרים כאלו) אלא אני בא להגיד שזיהוי דוגמאות בדאטהסט הדומות מדי לבנצ''מרקים אינן מצליחות לפלטר את הדוגמאות האלו. והתוצאה היא מודל שהוא אוברפיט על בנצמרק כזה או אחר.'

# --- Paragraph 4 ---
$p4 = $d.Paragraphs.Item(4).Range
$found4 = $p4.Find.Execute('את שרשראות ה-reasoning לשאלות שהצליחו להגיע לתשובה נכונה מוסיפים לסט שנקרא לו D_N. לשאלות שהמודל לא הצליח להגיע לתשובה סופית נכונה אנחנו מוסיפים רמז (במאמר זה נקרא rationalization) שעוזר למודל לבנות את שרשרת ה-reasoning. השאלות שהצליחו להגיע לתשובה הנכונה אחרי הרמז גם נוספים ל D_N. לאחר מכן מבצעים איטרציה אחת של שיטת מורד הגרדיאנט נבחרת על D_N ומעדכנים את משקלי המודל. חוזרים על השלבים האלו עד שהלוס מתייצב.', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Assert-True $found4 "paragraph 4 old text not found"
$p4.Text = 'כאמור יש שיטות די בסיסיות הבודקות את הדמיון בין הדוגמאות בדאטהסט לדוגמאות בבנצ''מארק מבוססות על n-grams ועל דמיון סמנטי המחושב באמצעות מרחק בין הייצוגי של הדוגמאות בדאטהסט ובבנצ''מרק. המאמר המסוקר טוען שזה לא מספיק וצריך לעשות בדיקה נוספת לזיהוי של דוגמאות אלו. בגדול המאמר מציע בנוסף לבדיקה הסמנטית לרתום איזה LLM עוצמתי לבדיקה של דמיון דוגמאות. '

# --- Paragraph 5 ---
$p5 = $d.Paragraphs.Item(5).Range
$found5 = $p5.Find.Execute('זהו זה, שיטה אינטואיטיבית ופשוטה שקיבלה כמה מאמרי השמך די כבדים שבתקווה אסקור אותם גם כן ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Assert-True $found5 "paragraph 5 old text not found"
$p5.Text = 'בגדול מזהים K דוגמאות הכי דומות סמנטית לכל דוגמא בבנצ''מרק ואז מפעילים LLM חזק כמו GPT4 עם איזה פרומפט מתוחכם כדי לזהות את הדוגמאות הבאמת דומות. המאמר מראה כי בצורה כזו הצליחו לתפוס דוגמאות שלמרות שנראות שונה מהוות rephrasing של דוגמא מסוימת מהבנצ''מרק. ואז מעיפים את הדוגמה הזו מהדאטהסט. '

# --- Paragraph 6 ---
$p6 = $d.Paragraphs.Item(6).Range
$found6 = $p6.Find.Execute('https://arxiv.org/pdf/2203.14465', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Assert-True $found6 "paragraph 6 old text not found"
$p6.Text = 'המאמר טוען כי ללא שימוש בשיטה שלהם ניתן ״לאמן״ מודל 13B כדי ש״ינצח״ את GPT4 על כל הבנצ''מרקים - נצחון לא אמיתי אמנם.'

# --- Append two brand-new paragraphs (style "Normal", inherited automatically) ---
# --- after paragraph 6, pushing the (already-updated) link paragraph to slot 8 ---
$p6r = $d.Paragraphs.Item(6).Range
$p6r.InsertParagraphAfter()
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Text = 'מאמר ללא יותר מדי חדשנות אך מעלה נושא מאד מעניין'

$p7r = $d.Paragraphs.Item(7).Range
$p7r.InsertParagraphAfter()
$p8 = $d.Paragraphs.Item(8)
$p8.Range.Text = 'https://arxiv.org/pdf/2311.04850'

Assert-True ($d.Paragraphs.Count -eq 8) "expected 8 paragraphs after edit"
Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
